# Applies the "Add files via upload" edit to the coordinates worksheet.
# Changes:
#  - New column G ("zdjecie") added; G2 gets a filenames value.
#  - Columns A (x), B (y) and F (pietro) get updated numeric data for rows 2-16.
#  - Row 17 (previously "Szumlas Emma") is cleared out entirely.
#  - The sheet view's top-left cell and active selection move back to A1 / G10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header + first data value
$ws.Range("G1").Value = "zdjecie"
$ws.Range("G2").Value = "1981-boze-cialo.jpg;1996_KEN-Budowa.jpg;KEN-Bartoka.jpg"

# Updated coordinate (A=x) / size (B=y) / floor (F=pietro) values per row
$data = @{
    2  = @{ A = 1000; B = 200;  F = 1 }
    3  = @{ A = 1000; B = 300;  F = 0 }
    4  = @{ A = 900;  B = 1500; F = 3 }
    5  = @{ A = 1435; B = 1400; F = 2 }
    6  = @{ A = 630;  B = 1500; F = 2 }
    7  = @{ A = 1200; B = 800;  F = 2 }
    8  = @{ A = 800;  B = 800;  F = 2 }
    9  = @{ A = 1390; B = 800;  F = 1 }
    10 = @{ A = 1500; B = 1500; F = 3 }
    11 = @{ A = 1030; B = 1450; F = 1 }
    12 = @{ A = 1435; B = 700;  F = 2 }
    13 = @{ A = 1250; B = 1400; F = 2 }
    14 = @{ A = 1550; B = 850;  F = 3 }
    15 = @{ A = 1100; B = 300;  F = 0 }
    16 = @{ A = 980;  B = 1500; F = 2 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("A$row").Value = $vals.A
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("F$row").Value = $vals.F
}

# Row 17 used to hold "Szumlas Emma" - it is now entirely empty
$ws.Range("A17:F17").ClearContents()

# Restore view: top-left cell back to A1, selection moved to G10
$excel.Goto($ws.Range("A1"))
$ws.Range("G10").Select()
